# Danh_sách_công_việc.xlsx - "chinh sua logic sap xep cong viec ca nhan"
# Appends two new personal-task rows (id 190 and 191) to the task list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The start_date/end_date columns (C, D) hold plain text such as "2021-04-13"
# rather than real Excel dates. Force the incoming cells to Text format first
# so Excel doesn't auto-convert these date-look-alike strings into date
# serial numbers; we strip the temporary formatting again afterwards so the
# cells end up stored exactly like the existing text-only date cells.
$ws.Range("C39:D40").NumberFormat = "@"

# Row 39: id=190
$ws.Range("A39").Value = 190
$ws.Range("B39").Value = "zxczx"
$ws.Range("C39").Value = "2021-04-13"
$ws.Range("D39").Value = "2021-04-17"
$ws.Range("E39").Value = "Chưa hoàn thành"

# Row 40: id=191
$ws.Range("A40").Value = 191
$ws.Range("B40").Value = "sczxc"
$ws.Range("C40").Value = "2021-04-25"
$ws.Range("D40").Value = "2021-04-27"
$ws.Range("E40").Value = "Chưa hoàn thành"

# Drop the temporary Text number format so the new cells match the
# General-formatted, string-valued cells used throughout the rest of the
# sheet.
$ws.Range("C39:D40").ClearFormats()
